$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(-83.0938, 490, 50.4013, 0.4044, -85.6212, 490, 50.4013, 0.3924, -27.6979, 490, 50.4013, 1.2131, 43392.3401, 113.3993, 51.4338, 6.6138)
    3  = @(-93.7182, 490, 50.4013, 0.3585, -92.878, 490, 50.4013, 0.3618, -31.2394, 490, 50.4013, 1.0756, 43392.3401, 113.3993, 51.4338, 6.6138)
    4  = @(-93.4259, 490, 50.4013, 0.3597, -93.151, 490, 50.4013, 0.3607, -31.142, 490, 50.4013, 1.079, 43392.3401, 113.3993, 51.4338, 6.6138)
    5  = @(-92.8758, 490, 50.4013, 0.3618, -92.8957, 490, 50.4013, 0.3617, -30.9586, 490, 50.4013, 1.0853, 43392.3401, 113.3993, 51.4338, 6.6138)
    6  = @(-92.6444, 490, 50.4013, 0.3627, -92.7045, 490, 50.4013, 0.3625, -30.8815, 490, 50.4013, 1.0881, 43392.3401, 113.3993, 51.4338, 6.6138)
    7  = @(-92.8758, 490, 50.4013, 0.3618, -92.8957, 490, 50.4013, 0.3617, -30.9586, 490, 50.4013, 1.0853, 43392.3401, 113.3993, 51.4338, 6.6138)
    8  = @(-93.4259, 490, 50.4013, 0.3597, -93.151, 490, 50.4013, 0.3607, -31.142, 490, 50.4013, 1.079, 43392.3401, 113.3993, 51.4338, 6.6138)
    9  = @(-93.7182, 490, 50.4013, 0.3585, -92.878, 490, 50.4013, 0.3618, -31.2394, 490, 50.4013, 1.0756, 43392.3401, 113.3993, 51.4338, 6.6138)
    10 = @(-83.0938, 490, 50.4013, 0.4044, -85.6212, 490, 50.4013, 0.3924, -27.6979, 490, 50.4013, 1.2131, 43392.3401, 113.3993, 51.4338, 6.6138)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
